# Updates cryptos.xlsx price/volume figures (and the EURNeutrino/BinanceUSD
# rank swap at rows 29-30) to match the "Wed Nov  1 23:50:36 UTC 2023" refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (D2, E2)
$ws.Range('D2').Value = '35.516.11'
$ws.Range('E2').Value = '  +2.47%  '

# Row 3 (D3, E3)
$ws.Range('D3').Value = '1.848.85'
$ws.Range('E3').Value = '  +1.90%  '

# Row 4 (E4)
$ws.Range('E4').Value = '  +0.07%  '

# Row 5 (D5, E5)
$ws.Range('D5').Value = '''228.21'
$ws.Range('E5').Value = '  +0.77%  '

# Row 6 (E6)
$ws.Range('E6').Value = '  +1.73%  '

# Row 7 (E7)
$ws.Range('E7').Value = '  +0.11%  '

# Row 8 (D8, E8)
$ws.Range('D8').Value = '''41.46'
$ws.Range('E8').Value = '  +7.70%  '

# Row 9 (D9, E9)
$ws.Range('D9').Value = '''0.308'
$ws.Range('E9').Value = '  +5.10%  '

# Row 10 (D10, E10)
$ws.Range('D10').Value = '''0.0689'
$ws.Range('E10').Value = '  +1.01%  '

# Row 11 (D11, E11)
$ws.Range('D11').Value = '''0.100'
$ws.Range('E11').Value = '  +3.27%  '

# Row 12 (D12, E12)
$ws.Range('D12').Value = '2.114.72'
$ws.Range('E12').Value = '  +1.85%  '

# Row 13 (D13, E13)
$ws.Range('D13').Value = '''11.64'
$ws.Range('E13').Value = '  +2.57%  '

# Row 14 (D14, E14)
$ws.Range('D14').Value = '1.846.00'
$ws.Range('E14').Value = '  +1.64%  '

# Row 15 (D15, E15)
$ws.Range('D15').Value = '''4.74'
$ws.Range('E15').Value = '  +6.39%  '

# Row 16 (D16, E16)
$ws.Range('D16').Value = '''0.668'
$ws.Range('E16').Value = '  +5.10%  '

# Row 17 (D17, E17)
$ws.Range('D17').Value = '35.455.64'
$ws.Range('E17').Value = '  +2.45%  '

# Row 18 (D18, E18)
$ws.Range('D18').Value = '''69.99'
$ws.Range('E18').Value = '  +1.68%  '

# Row 19 (D19, E19)
$ws.Range('D19').Value = '''245.70'
$ws.Range('E19').Value = '  +0.23%  '

# Row 20 (D20, E20)
$ws.Range('D20').Value = '0.0₃0795'
$ws.Range('E20').Value = '  +2.23%  '

# Row 21 (D21, E21)
$ws.Range('D21').Value = '''12.21'
$ws.Range('E21').Value = '  +7.89%  '

# Row 22 (D22, E22)
$ws.Range('D22').Value = '''4.76'
$ws.Range('E22').Value = '  +14.92%  '

# Row 23 (E23)
$ws.Range('E23').Value = '  +0.17%  '

# Row 24 (D24, E24)
$ws.Range('D24').Value = '''2.21'
$ws.Range('E24').Value = '  -0.81%  '

# Row 25 (D25, E25)
$ws.Range('D25').Value = '''171.57'
$ws.Range('E25').Value = '  -0.45%  '

# Row 26 (D26, E26)
$ws.Range('D26').Value = '''7.92'
$ws.Range('E26').Value = '  -0.01%  '

# Row 27 (D27)
$ws.Range('D27').Value = '''17.86'

# Row 28 (E28)
$ws.Range('E28').Value = '  +0.97%  '

# Row 29 (B29, C29, D29, E29)
$ws.Range('B29').Value = 'BinanceUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.14%  '

# Row 30 (B30, C30, D30, E30)
$ws.Range('B30').Value = 'EURNeutrino'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7YKHKSdb+eurneutrino-eurn'
$ws.Range('D30').Value = '3.409.53'
$ws.Range('E30').Value = '  +40.33%  '

# Row 31 (E31)
$ws.Range('E31').Value = '  +7.62%  '

# Row 32 (D32, E32)
$ws.Range('D32').Value = '''3.94'
$ws.Range('E32').Value = '  +3.24%  '

# Row 33 (D33, E33)
$ws.Range('D33').Value = '''4.06'
$ws.Range('E33').Value = '  +2.99%  '

# Row 34 (D34, E34)
$ws.Range('D34').Value = '''0.0536'
$ws.Range('E34').Value = '  +2.18%  '

# Row 35 (E35)
$ws.Range('E35').Value = '  +3.21%  '

# Row 36 (D36, E36)
$ws.Range('D36').Value = '''0.679'
$ws.Range('E36').Value = '  +3.25%  '

# Row 37 (D37, E37)
$ws.Range('D37').Value = '''1.04'
$ws.Range('E37').Value = '  +9.75%  '

# Row 38 (D38, E38)
$ws.Range('D38').Value = '''89.09'
$ws.Range('E38').Value = '  +9.56%  '

# Row 39 (D39, E39)
$ws.Range('D39').Value = '1.342.36'
$ws.Range('E39').Value = '  -1.93%  '

# Row 40 (E40)
$ws.Range('E40').Value = '  +1.43%  '

# Row 41 (D41, E41)
$ws.Range('D41').Value = '''0.0195'
$ws.Range('E41').Value = '  +3.26%  '

# Row 42 (E42)
$ws.Range('E42').Value = '  +1.59%  '

# Row 43 (E43)
$ws.Range('E43').Value = '  +3.85%  '

# Row 44 (D44, E44)
$ws.Range('D44').Value = '''14.91'
$ws.Range('E44').Value = '  +5.50%  '

# Row 45 (E45)
$ws.Range('E45').Value = '  +1.02%  '

# Row 46 (E46)
$ws.Range('E46').Value = '  +1.40%  '

# Row 47 (E47)
$ws.Range('E47').Value = '  +3.85%  '

# Row 48 (D48, E48)
$ws.Range('D48').Value = '''6.05'
$ws.Range('E48').Value = '  +4.53%  '

# Row 49 (D49, E49)
$ws.Range('D49').Value = '2.014.55'
$ws.Range('E49').Value = '  +1.90%  '

# Row 50 (D50, E50)
$ws.Range('D50').Value = '''104.65'
$ws.Range('E50').Value = '  +1.45%  '

# Row 51 (E51)
$ws.Range('E51').Value = '  +0.10%  '
